$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 128, pushing existing rows 128:165 down to 129:166.
$ws.Range("A128:R128").EntireRow.Insert()

# Populate the newly inserted row 128 with the new record.
$ws.Range("A128").Value = 4
$ws.Range("B128").Value = 'Feria Lagunitas de Puerto Montt'
$ws.Range("C128").Value = 'Los Lagos'
$ws.Range("D128").Value = 44663
$ws.Range("E128").Value = 10
$ws.Range("F128").Value = 100112009
$ws.Range("G128").Value = 'Acelga'
$ws.Range("H128").Value = 'Sin especificar'
$ws.Range("I128").Value = 'Primera'
$ws.Range("J128").Value = 80
$ws.Range("K128").Value = 10000
$ws.Range("L128").Value = 10000
$ws.Range("M128").Value = 10000
$ws.Range("N128").Value = '$/docena de atados (12 kilos)'
$ws.Range("O128").Value = 'Región de La Araucanía'
$ws.Range("P128").Value = 833
$ws.Range("Q128").Value = 12
$ws.Range("R128").Value = 'Hortaliza'

# Match the date-number-format style used by the rest of column D.
$ws.Range("D128").NumberFormat = $ws.Range("D129").NumberFormat
